$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 13-35: dataset correction (titolo/giornale/social/hate_speech realigned)
$ws.Cells.Item(13, 2).Value = 'L''implosione del sottomarino Titan'
$ws.Cells.Item(13, 3).Value = 'FanPage'
$ws.Cells.Item(13, 4).Value = 'Instagram'
$ws.Cells.Item(13, 5).Value = 3

$ws.Cells.Item(14, 2).Value = 'L''implosione del sottomarino Titan'
$ws.Cells.Item(14, 3).Value = 'FanPage'
$ws.Cells.Item(14, 4).Value = 'YouTube'
$ws.Cells.Item(14, 5).Value = 3

$ws.Cells.Item(15, 2).Value = 'L''implosione del sottomarino Titan'
$ws.Cells.Item(15, 3).Value = 'Il Corriere Della Sera'
$ws.Cells.Item(15, 4).Value = 'Facebook'
$ws.Cells.Item(15, 5).Value = 2

$ws.Cells.Item(16, 2).Value = 'L''implosione del sottomarino Titan'
$ws.Cells.Item(16, 3).Value = 'Il Corriere Della Sera'
$ws.Cells.Item(16, 4).Value = 'Instagram'
$ws.Cells.Item(16, 5).Value = 4

$ws.Cells.Item(17, 2).Value = 'L''implosione del sottomarino Titan'
$ws.Cells.Item(17, 3).Value = 'Il Corriere Della Sera'
$ws.Cells.Item(17, 4).Value = 'YouTube'
$ws.Cells.Item(17, 5).Value = 1

$ws.Cells.Item(18, 2).Value = 'L''implosione del sottomarino Titan'
$ws.Cells.Item(18, 3).Value = 'La Repubblica'
$ws.Cells.Item(18, 4).Value = 'Facebook'
$ws.Cells.Item(18, 5).Value = 1

$ws.Cells.Item(19, 2).Value = 'L''implosione del sottomarino Titan'
$ws.Cells.Item(19, 3).Value = 'La Repubblica'
$ws.Cells.Item(19, 4).Value = 'Instagram'
$ws.Cells.Item(19, 5).Value = 10

$ws.Cells.Item(20, 2).Value = 'L''implosione del sottomarino Titan'
$ws.Cells.Item(20, 3).Value = 'La Repubblica'
$ws.Cells.Item(20, 4).Value = 'YouTube'
$ws.Cells.Item(20, 5).Value = 6

$ws.Cells.Item(21, 2).Value = 'L''omicidio di Giulia Cecchettin'
$ws.Cells.Item(21, 3).Value = 'FanPage'
$ws.Cells.Item(21, 4).Value = 'Facebook'
$ws.Cells.Item(21, 5).Value = 7

$ws.Cells.Item(22, 2).Value = 'L''omicidio di Giulia Cecchettin'
$ws.Cells.Item(22, 3).Value = 'FanPage'
$ws.Cells.Item(22, 4).Value = 'Instagram'
$ws.Cells.Item(22, 5).Value = 4

$ws.Cells.Item(23, 2).Value = 'L''omicidio di Giulia Cecchettin'
$ws.Cells.Item(23, 3).Value = 'FanPage'
$ws.Cells.Item(23, 4).Value = 'YouTube'
$ws.Cells.Item(23, 5).Value = 1

$ws.Cells.Item(24, 2).Value = 'L''omicidio di Giulia Cecchettin'
$ws.Cells.Item(24, 3).Value = 'Il Corriere Della Sera'
$ws.Cells.Item(24, 4).Value = 'Facebook'
$ws.Cells.Item(24, 5).Value = 6

$ws.Cells.Item(25, 2).Value = 'L''omicidio di Giulia Cecchettin'
$ws.Cells.Item(25, 3).Value = 'Il Corriere Della Sera'
$ws.Cells.Item(25, 4).Value = 'Instagram'
$ws.Cells.Item(25, 5).Value = 10

$ws.Cells.Item(26, 2).Value = 'L''omicidio di Giulia Cecchettin'
$ws.Cells.Item(26, 3).Value = 'Il Corriere Della Sera'
$ws.Cells.Item(26, 4).Value = 'YouTube'
$ws.Cells.Item(26, 5).Value = 7

$ws.Cells.Item(27, 2).Value = 'L''omicidio di Giulia Cecchettin'
$ws.Cells.Item(27, 3).Value = 'La Repubblica'
$ws.Cells.Item(27, 4).Value = 'Facebook'
$ws.Cells.Item(27, 5).Value = 16

$ws.Cells.Item(28, 2).Value = 'L''omicidio di Giulia Cecchettin'
$ws.Cells.Item(28, 3).Value = 'La Repubblica'
$ws.Cells.Item(28, 4).Value = 'Instagram'
$ws.Cells.Item(28, 5).Value = 16

$ws.Cells.Item(29, 2).Value = 'L''omicidio di Giulia Cecchettin'
$ws.Cells.Item(29, 3).Value = 'La Repubblica'
$ws.Cells.Item(29, 4).Value = 'YouTube'
$ws.Cells.Item(29, 5).Value = 6

$ws.Cells.Item(30, 2).Value = 'Strage di Cutro'
$ws.Cells.Item(30, 3).Value = 'FanPage'
$ws.Cells.Item(30, 4).Value = 'Facebook'
$ws.Cells.Item(30, 5).Value = 9

$ws.Cells.Item(31, 2).Value = 'Strage di Cutro'
$ws.Cells.Item(31, 3).Value = 'FanPage'
$ws.Cells.Item(31, 4).Value = 'Instagram'
$ws.Cells.Item(31, 5).Value = 32

$ws.Cells.Item(32, 2).Value = 'Strage di Cutro'
$ws.Cells.Item(32, 3).Value = 'FanPage'
$ws.Cells.Item(32, 4).Value = 'YouTube'
$ws.Cells.Item(32, 5).Value = 13

$ws.Cells.Item(33, 2).Value = 'Strage di Cutro'
$ws.Cells.Item(33, 3).Value = 'Il Corriere Della Sera'
$ws.Cells.Item(33, 4).Value = 'Facebook'
$ws.Cells.Item(33, 5).Value = 28

$ws.Cells.Item(34, 2).Value = 'Strage di Cutro'
$ws.Cells.Item(34, 3).Value = 'Il Corriere Della Sera'
$ws.Cells.Item(34, 4).Value = 'Instagram'
$ws.Cells.Item(34, 5).Value = 15

$ws.Cells.Item(35, 2).Value = 'Strage di Cutro'
$ws.Cells.Item(35, 3).Value = 'Il Corriere Della Sera'
$ws.Cells.Item(35, 4).Value = 'YouTube'
$ws.Cells.Item(35, 5).Value = 10

# Rows 36-38: new rows appended for "Strage di Cutro" / La Repubblica
$aCell = $ws.Cells.Item(36, 1)
$aCell.Value = "'"
$aCell.Style = 'Normal'
$ws.Cells.Item(36, 2).Value = 'Strage di Cutro'
$ws.Cells.Item(36, 3).Value = 'La Repubblica'
$ws.Cells.Item(36, 4).Value = 'Facebook'
$ws.Cells.Item(36, 5).Value = 39

$aCell = $ws.Cells.Item(37, 1)
$aCell.Value = "'"
$aCell.Style = 'Normal'
$ws.Cells.Item(37, 2).Value = 'Strage di Cutro'
$ws.Cells.Item(37, 3).Value = 'La Repubblica'
$ws.Cells.Item(37, 4).Value = 'Instagram'
$ws.Cells.Item(37, 5).Value = 13

$aCell = $ws.Cells.Item(38, 1)
$aCell.Value = "'"
$aCell.Style = 'Normal'
$ws.Cells.Item(38, 2).Value = 'Strage di Cutro'
$ws.Cells.Item(38, 3).Value = 'La Repubblica'
$ws.Cells.Item(38, 4).Value = 'YouTube'
$ws.Cells.Item(38, 5).Value = 14

Write-Output "edit applied"